$d = $word.ActiveDocument

# Helper: replace the whole text of $range with $newText using Find/Replace
# scoped to that range (Range.Text = "..." in this host inserts rather than
# overwrites, so Find/Execute is used for reliable in-place replacement).
function Set-RangeText($range, $oldText, $newText) {
    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# --- 1. Title / headline text: replace both occurrences --------------------
# (Heading1 at top of doc, and the bold run near the bottom)
$d.Content.Find.Execute(
    "Play Castle of Terror for Free - Review of Big Time Gaming's Horror-themed Slot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Castle of Terror Free - Review of Gameplay, Graphics, and Bonus Features",
    2) | Out-Null

# --- 2. "What we like" bullet list ------------------------------------------
# Before:
#   High win potential
#   Impressive bonus features
#   Exceptional graphics and gameplay
#   Mobile-optimized
# After:
#   Exceptional graphics and gameplay
#   Well-crafted horror theme
#   High win potential and volatility
#   Impressive bonus features

# Find the bullet paragraph that currently reads "High win potential" and
# turn it into the new first bullet, then insert the two new bullets right
# after it (inheriting the ListBullet paragraph style automatically).
# NOTE: Paragraph.Range.Text includes the trailing paragraph mark (vbCr), so
# trim it off before comparing against plain literal strings.
$idxHigh = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "High win potential") {
        $idxHigh = $i
        break
    }
}

$pHighRange = $d.Paragraphs.Item($idxHigh).Range
Set-RangeText $pHighRange "High win potential" "Exceptional graphics and gameplay"

$d.Paragraphs.Item($idxHigh).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idxHigh + 1).Range.InsertBefore("Well-crafted horror theme")

$d.Paragraphs.Item($idxHigh + 1).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($idxHigh + 2).Range.InsertBefore("High win potential and volatility")

# "Impressive bonus features" bullet stays as-is (no change needed).

# Remove the now-duplicate "Exceptional graphics and gameplay" bullet and the
# "Mobile-optimized" bullet that used to follow "Impressive bonus features".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Mobile-optimized") {
        $d.Paragraphs.Item($i).Range.Delete() | Out-Null
        break
    }
}
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if (($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Exceptional graphics and gameplay") -and ($i -ne $idxHigh)) {
        $d.Paragraphs.Item($i).Range.Delete() | Out-Null
        break
    }
}

# --- 3. "What we don't like" bullet list ------------------------------------
# Before:
#   High volatility
#   Limited betting range
# After:
#   Limited betting range
#   Standard RTP rate
# Process the second bullet first so the still-unique old text can be found
# before the first bullet is renamed to the same text as the (old) second one.
$d.Content.Find.Execute("Limited betting range", $true, $false, $false, $false, $false, $true, 1, $false, "Standard RTP rate", 2) | Out-Null
$d.Content.Find.Execute("High volatility", $true, $false, $false, $false, $false, $true, 1, $false, "Limited betting range", 2) | Out-Null

# --- 4. Meta/subtitle italic text -------------------------------------------
$d.Content.Find.Execute(
    "Experience terror with Castle of Terror, a horror-themed online slot game by Big Time Gaming. Play for free and read our review here to learn more.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Castle of Terror, a horror-themed slot game with exceptional graphics and impressive bonus features. Play for free and experience the terror!",
    2) | Out-Null
